# Update the "Backlog" sheet to mark rows 104-106 and 109 as complete (column C = "X"),
# and move the sheet's view/selection down to reflect the new rows added at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

# Mark the newly-finished story tasks with an "X" in column C, matching the
# existing convention used throughout the sheet (see rows 79-103, 107-108, etc.).
$ws.Range("C104").Value = "X"
$ws.Range("C105").Value = "X"
$ws.Range("C106").Value = "X"
$ws.Range("C109").Value = "X"

# Recalculate so the COUNTA summary in D109 picks up the new entries.
$excel.CalculateFull()

# Scroll the view down and move the active selection to reflect where the
# user was last working in the sheet.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 86
$win.ScrollColumn = 1
$ws.Range("C110").Select()
